# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect freshly-scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsAll  = $wb.Worksheets.Item("全部类型")

# Rows that line up 1:1 between the two sheets (same row numbers, simple +1 bump)
$wsExpo.Range("F3").Value  = 20326
$wsExpo.Range("F4").Value  = 804
$wsExpo.Range("F8").Value  = 7648
$wsExpo.Range("F11").Value = 279
$wsExpo.Range("F13").Value = 165
$wsExpo.Range("F15").Value = 18
$wsExpo.Range("F19").Value = 452
$wsExpo.Range("F26").Value = 1123
$wsExpo.Range("F33").Value = 4002
$wsExpo.Range("F37").Value = 12733

$wsAll.Range("F3").Value  = 20326
$wsAll.Range("F4").Value  = 804
$wsAll.Range("F8").Value  = 7648
$wsAll.Range("F11").Value = 279
$wsAll.Range("F13").Value = 165
$wsAll.Range("F15").Value = 18
$wsAll.Range("F19").Value = 452
$wsAll.Range("F26").Value = 1123
$wsAll.Range("F36").Value = 4004
$wsAll.Range("F40").Value = 12733
